$wb = $excel.ActiveWorkbook

# Rename sheets:
#   "principal" -> "tipo_persona"
#   "relacion"  -> "tipo_persona_rel"
$wb.Worksheets.Item("principal").Name = "tipo_persona"
$wb.Worksheets.Item("relacion").Name = "tipo_persona_rel"
